# Update odds values on Sheet1 to match the latest FlashScore scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Ind. Medellin vs Alianza)
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 9
$ws.Range("O2").Value = 1.4
$ws.Range("P2").Value = 2.75
$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.65

# Row 4 (Correcaminos vs Atl. Morelia)
$ws.Range("K4").Value = 2.12
$ws.Range("L4").Value = 2.65
$ws.Range("AR4").Value = 110
$ws.Range("AT4").Value = 2.87
$ws.Range("AZ4").Value = 40
